$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary cell updates ---
$ws.Range("E11").Value = 4785906
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 109

# --- Move the footer ("signature") block down from rows 128-129 to rows 133-134 FIRST,
# before row 128/123 formatting is touched, so the copied borders are the original ones. ---
$ws.Range("B128:C129").Copy()
$ws.Range("B133:C134").PasteSpecial(-4122)
$ws.Range("H128:J129").Copy()
$ws.Range("H133:J134").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Footer text content at the new location (values copied along are overwritten to be safe).
$ws.Range("B133").Value = "___________________________________"
$ws.Range("H133").Value = "___________________________________"
$ws.Range("B134").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H134").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# Re-merge the footer cells at their new location.
$ws.Range("B133:C133").Merge()
$ws.Range("H133:J133").Merge()
$ws.Range("B134:C134").Merge()
$ws.Range("H134:J134").Merge()

# Drop the old footer merges and clear rows 128-129 of the old footer content/format.
$ws.Range("B128:C128").UnMerge()
$ws.Range("H128:J128").UnMerge()
$ws.Range("B129:C129").UnMerge()
$ws.Range("H129:J129").UnMerge()
$ws.Range("B129:J129").Clear()

# --- Prepare formatting for the data block ---
# Copy the "last row" format (currently on row 123) down to row 128, which will
# become the new final data row of the table.
$ws.Range("B123:J123").Copy()
$ws.Range("B128:J128").PasteSpecial(-4122)

# Copy the "normal row" format (row 122) onto row 123 and the newly used rows 124-127.
$ws.Range("B122:J122").Copy()
$ws.Range("B123:J127").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Rewrite the account-statement data table ---
# Row 16-124: worker CC 10236900 OLMEDO DE JESUS CASTRO CAMPEON, periods 2507..1607 descending.
# Row 125-127: worker CC 73561364 RAFAEL ENRIQUE JIMENEZ GONGORA.
# Row 128: worker CC 8802936 CARLOS ALBERTO ARZUZA DIAZ (final row, "last" style).
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "10236900"
$ws.Range("D16").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 43480
$ws.Range("G16").Value = 1087000
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "10236900"
$ws.Range("D17").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 43480
$ws.Range("G17").Value = 1087000
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "10236900"
$ws.Range("D18").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 43480
$ws.Range("G18").Value = 1087000
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "10236900"
$ws.Range("D19").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 43480
$ws.Range("G19").Value = 1087000
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "10236900"
$ws.Range("D20").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E20").Value = "2503"
$ws.Range("F20").Value = 43480
$ws.Range("G20").Value = 1087000
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "10236900"
$ws.Range("D21").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E21").Value = "2502"
$ws.Range("F21").Value = 43480
$ws.Range("G21").Value = 1087000
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "10236900"
$ws.Range("D22").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E22").Value = "2501"
$ws.Range("F22").Value = 43480
$ws.Range("G22").Value = 1087000
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "10236900"
$ws.Range("D23").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E23").Value = "2412"
$ws.Range("F23").Value = 43480
$ws.Range("G23").Value = 1087000
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "10236900"
$ws.Range("D24").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E24").Value = "2411"
$ws.Range("F24").Value = 43480
$ws.Range("G24").Value = 1087000
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "10236900"
$ws.Range("D25").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E25").Value = "2410"
$ws.Range("F25").Value = 43480
$ws.Range("G25").Value = 1087000
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "10236900"
$ws.Range("D26").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E26").Value = "2409"
$ws.Range("F26").Value = 43480
$ws.Range("G26").Value = 1087000
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "10236900"
$ws.Range("D27").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E27").Value = "2408"
$ws.Range("F27").Value = 43480
$ws.Range("G27").Value = 1087000
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "10236900"
$ws.Range("D28").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E28").Value = "2407"
$ws.Range("F28").Value = 43480
$ws.Range("G28").Value = 1087000
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "10236900"
$ws.Range("D29").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E29").Value = "2406"
$ws.Range("F29").Value = 43480
$ws.Range("G29").Value = 1087000
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "10236900"
$ws.Range("D30").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E30").Value = "2405"
$ws.Range("F30").Value = 43480
$ws.Range("G30").Value = 1087000
$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "10236900"
$ws.Range("D31").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E31").Value = "2404"
$ws.Range("F31").Value = 43480
$ws.Range("G31").Value = 1087000
$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "10236900"
$ws.Range("D32").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E32").Value = "2403"
$ws.Range("F32").Value = 43480
$ws.Range("G32").Value = 1087000
$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "10236900"
$ws.Range("D33").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E33").Value = "2402"
$ws.Range("F33").Value = 43480
$ws.Range("G33").Value = 1087000
$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "10236900"
$ws.Range("D34").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E34").Value = "2401"
$ws.Range("F34").Value = 43480
$ws.Range("G34").Value = 1087000
$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "10236900"
$ws.Range("D35").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E35").Value = "2312"
$ws.Range("F35").Value = 43480
$ws.Range("G35").Value = 1087000
$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "10236900"
$ws.Range("D36").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E36").Value = "2311"
$ws.Range("F36").Value = 43480
$ws.Range("G36").Value = 1087000
$ws.Range("B37").Value = "CC"
$ws.Range("C37").Value = "10236900"
$ws.Range("D37").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E37").Value = "2310"
$ws.Range("F37").Value = 43480
$ws.Range("G37").Value = 1087000
$ws.Range("B38").Value = "CC"
$ws.Range("C38").Value = "10236900"
$ws.Range("D38").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E38").Value = "2309"
$ws.Range("F38").Value = 43480
$ws.Range("G38").Value = 1087000
$ws.Range("B39").Value = "CC"
$ws.Range("C39").Value = "10236900"
$ws.Range("D39").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E39").Value = "2308"
$ws.Range("F39").Value = 43480
$ws.Range("G39").Value = 1087000
$ws.Range("B40").Value = "CC"
$ws.Range("C40").Value = "10236900"
$ws.Range("D40").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E40").Value = "2307"
$ws.Range("F40").Value = 43480
$ws.Range("G40").Value = 1087000
$ws.Range("B41").Value = "CC"
$ws.Range("C41").Value = "10236900"
$ws.Range("D41").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E41").Value = "2306"
$ws.Range("F41").Value = 43480
$ws.Range("G41").Value = 1087000
$ws.Range("B42").Value = "CC"
$ws.Range("C42").Value = "10236900"
$ws.Range("D42").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E42").Value = "2305"
$ws.Range("F42").Value = 43480
$ws.Range("G42").Value = 1087000
$ws.Range("B43").Value = "CC"
$ws.Range("C43").Value = "10236900"
$ws.Range("D43").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E43").Value = "2304"
$ws.Range("F43").Value = 43480
$ws.Range("G43").Value = 1087000
$ws.Range("B44").Value = "CC"
$ws.Range("C44").Value = "10236900"
$ws.Range("D44").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E44").Value = "2303"
$ws.Range("F44").Value = 43480
$ws.Range("G44").Value = 1087000
$ws.Range("B45").Value = "CC"
$ws.Range("C45").Value = "10236900"
$ws.Range("D45").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E45").Value = "2302"
$ws.Range("F45").Value = 43480
$ws.Range("G45").Value = 1087000
$ws.Range("B46").Value = "CC"
$ws.Range("C46").Value = "10236900"
$ws.Range("D46").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E46").Value = "2301"
$ws.Range("F46").Value = 43480
$ws.Range("G46").Value = 1087000
$ws.Range("B47").Value = "CC"
$ws.Range("C47").Value = "10236900"
$ws.Range("D47").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E47").Value = "2212"
$ws.Range("F47").Value = 43480
$ws.Range("G47").Value = 1087000
$ws.Range("B48").Value = "CC"
$ws.Range("C48").Value = "10236900"
$ws.Range("D48").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E48").Value = "2211"
$ws.Range("F48").Value = 43480
$ws.Range("G48").Value = 1087000
$ws.Range("B49").Value = "CC"
$ws.Range("C49").Value = "10236900"
$ws.Range("D49").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E49").Value = "2210"
$ws.Range("F49").Value = 43480
$ws.Range("G49").Value = 1087000
$ws.Range("B50").Value = "CC"
$ws.Range("C50").Value = "10236900"
$ws.Range("D50").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E50").Value = "2209"
$ws.Range("F50").Value = 43480
$ws.Range("G50").Value = 1087000
$ws.Range("B51").Value = "CC"
$ws.Range("C51").Value = "10236900"
$ws.Range("D51").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E51").Value = "2208"
$ws.Range("F51").Value = 43480
$ws.Range("G51").Value = 1087000
$ws.Range("B52").Value = "CC"
$ws.Range("C52").Value = "10236900"
$ws.Range("D52").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E52").Value = "2207"
$ws.Range("F52").Value = 43480
$ws.Range("G52").Value = 1087000
$ws.Range("B53").Value = "CC"
$ws.Range("C53").Value = "10236900"
$ws.Range("D53").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E53").Value = "2206"
$ws.Range("F53").Value = 43480
$ws.Range("G53").Value = 1087000
$ws.Range("B54").Value = "CC"
$ws.Range("C54").Value = "10236900"
$ws.Range("D54").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E54").Value = "2205"
$ws.Range("F54").Value = 43480
$ws.Range("G54").Value = 1087000
$ws.Range("B55").Value = "CC"
$ws.Range("C55").Value = "10236900"
$ws.Range("D55").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E55").Value = "2204"
$ws.Range("F55").Value = 43480
$ws.Range("G55").Value = 1087000
$ws.Range("B56").Value = "CC"
$ws.Range("C56").Value = "10236900"
$ws.Range("D56").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E56").Value = "2203"
$ws.Range("F56").Value = 43480
$ws.Range("G56").Value = 1087000
$ws.Range("B57").Value = "CC"
$ws.Range("C57").Value = "10236900"
$ws.Range("D57").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E57").Value = "2202"
$ws.Range("F57").Value = 43480
$ws.Range("G57").Value = 1087000
$ws.Range("B58").Value = "CC"
$ws.Range("C58").Value = "10236900"
$ws.Range("D58").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E58").Value = "2201"
$ws.Range("F58").Value = 43480
$ws.Range("G58").Value = 1087000
$ws.Range("B59").Value = "CC"
$ws.Range("C59").Value = "10236900"
$ws.Range("D59").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E59").Value = "2112"
$ws.Range("F59").Value = 43480
$ws.Range("G59").Value = 1087000
$ws.Range("B60").Value = "CC"
$ws.Range("C60").Value = "10236900"
$ws.Range("D60").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E60").Value = "2111"
$ws.Range("F60").Value = 43480
$ws.Range("G60").Value = 1087000
$ws.Range("B61").Value = "CC"
$ws.Range("C61").Value = "10236900"
$ws.Range("D61").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E61").Value = "2110"
$ws.Range("F61").Value = 43480
$ws.Range("G61").Value = 1087000
$ws.Range("B62").Value = "CC"
$ws.Range("C62").Value = "10236900"
$ws.Range("D62").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E62").Value = "2109"
$ws.Range("F62").Value = 43480
$ws.Range("G62").Value = 1087000
$ws.Range("B63").Value = "CC"
$ws.Range("C63").Value = "10236900"
$ws.Range("D63").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E63").Value = "2108"
$ws.Range("F63").Value = 43480
$ws.Range("G63").Value = 1087000
$ws.Range("B64").Value = "CC"
$ws.Range("C64").Value = "10236900"
$ws.Range("D64").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E64").Value = "2107"
$ws.Range("F64").Value = 43480
$ws.Range("G64").Value = 1087000
$ws.Range("B65").Value = "CC"
$ws.Range("C65").Value = "10236900"
$ws.Range("D65").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E65").Value = "2106"
$ws.Range("F65").Value = 43480
$ws.Range("G65").Value = 1087000
$ws.Range("B66").Value = "CC"
$ws.Range("C66").Value = "10236900"
$ws.Range("D66").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E66").Value = "2105"
$ws.Range("F66").Value = 43480
$ws.Range("G66").Value = 1087000
$ws.Range("B67").Value = "CC"
$ws.Range("C67").Value = "10236900"
$ws.Range("D67").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E67").Value = "2104"
$ws.Range("F67").Value = 43480
$ws.Range("G67").Value = 1087000
$ws.Range("B68").Value = "CC"
$ws.Range("C68").Value = "10236900"
$ws.Range("D68").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E68").Value = "2103"
$ws.Range("F68").Value = 43480
$ws.Range("G68").Value = 1087000
$ws.Range("B69").Value = "CC"
$ws.Range("C69").Value = "10236900"
$ws.Range("D69").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E69").Value = "2102"
$ws.Range("F69").Value = 43480
$ws.Range("G69").Value = 1087000
$ws.Range("B70").Value = "CC"
$ws.Range("C70").Value = "10236900"
$ws.Range("D70").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E70").Value = "2101"
$ws.Range("F70").Value = 43480
$ws.Range("G70").Value = 1087000
$ws.Range("B71").Value = "CC"
$ws.Range("C71").Value = "10236900"
$ws.Range("D71").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E71").Value = "2012"
$ws.Range("F71").Value = 43480
$ws.Range("G71").Value = 1087000
$ws.Range("B72").Value = "CC"
$ws.Range("C72").Value = "10236900"
$ws.Range("D72").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E72").Value = "2011"
$ws.Range("F72").Value = 43480
$ws.Range("G72").Value = 1087000
$ws.Range("B73").Value = "CC"
$ws.Range("C73").Value = "10236900"
$ws.Range("D73").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E73").Value = "2010"
$ws.Range("F73").Value = 43480
$ws.Range("G73").Value = 1087000
$ws.Range("B74").Value = "CC"
$ws.Range("C74").Value = "10236900"
$ws.Range("D74").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E74").Value = "2009"
$ws.Range("F74").Value = 43480
$ws.Range("G74").Value = 1087000
$ws.Range("B75").Value = "CC"
$ws.Range("C75").Value = "10236900"
$ws.Range("D75").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E75").Value = "2008"
$ws.Range("F75").Value = 43480
$ws.Range("G75").Value = 1087000
$ws.Range("B76").Value = "CC"
$ws.Range("C76").Value = "10236900"
$ws.Range("D76").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E76").Value = "2007"
$ws.Range("F76").Value = 43480
$ws.Range("G76").Value = 1087000
$ws.Range("B77").Value = "CC"
$ws.Range("C77").Value = "10236900"
$ws.Range("D77").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E77").Value = "2006"
$ws.Range("F77").Value = 43480
$ws.Range("G77").Value = 1087000
$ws.Range("B78").Value = "CC"
$ws.Range("C78").Value = "10236900"
$ws.Range("D78").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E78").Value = "2005"
$ws.Range("F78").Value = 43480
$ws.Range("G78").Value = 1087000
$ws.Range("B79").Value = "CC"
$ws.Range("C79").Value = "10236900"
$ws.Range("D79").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E79").Value = "2004"
$ws.Range("F79").Value = 43480
$ws.Range("G79").Value = 1087000
$ws.Range("B80").Value = "CC"
$ws.Range("C80").Value = "10236900"
$ws.Range("D80").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E80").Value = "2003"
$ws.Range("F80").Value = 43480
$ws.Range("G80").Value = 1087000
$ws.Range("B81").Value = "CC"
$ws.Range("C81").Value = "10236900"
$ws.Range("D81").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E81").Value = "2002"
$ws.Range("F81").Value = 43480
$ws.Range("G81").Value = 1087000
$ws.Range("B82").Value = "CC"
$ws.Range("C82").Value = "10236900"
$ws.Range("D82").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E82").Value = "2001"
$ws.Range("F82").Value = 43480
$ws.Range("G82").Value = 1087000
$ws.Range("B83").Value = "CC"
$ws.Range("C83").Value = "10236900"
$ws.Range("D83").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E83").Value = "1912"
$ws.Range("F83").Value = 43480
$ws.Range("G83").Value = 1087000
$ws.Range("B84").Value = "CC"
$ws.Range("C84").Value = "10236900"
$ws.Range("D84").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E84").Value = "1911"
$ws.Range("F84").Value = 43480
$ws.Range("G84").Value = 1087000
$ws.Range("B85").Value = "CC"
$ws.Range("C85").Value = "10236900"
$ws.Range("D85").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E85").Value = "1910"
$ws.Range("F85").Value = 43480
$ws.Range("G85").Value = 1087000
$ws.Range("B86").Value = "CC"
$ws.Range("C86").Value = "10236900"
$ws.Range("D86").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E86").Value = "1909"
$ws.Range("F86").Value = 43480
$ws.Range("G86").Value = 1087000
$ws.Range("B87").Value = "CC"
$ws.Range("C87").Value = "10236900"
$ws.Range("D87").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E87").Value = "1908"
$ws.Range("F87").Value = 43480
$ws.Range("G87").Value = 1087000
$ws.Range("B88").Value = "CC"
$ws.Range("C88").Value = "10236900"
$ws.Range("D88").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E88").Value = "1907"
$ws.Range("F88").Value = 43480
$ws.Range("G88").Value = 1087000
$ws.Range("B89").Value = "CC"
$ws.Range("C89").Value = "10236900"
$ws.Range("D89").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E89").Value = "1906"
$ws.Range("F89").Value = 43480
$ws.Range("G89").Value = 1087000
$ws.Range("B90").Value = "CC"
$ws.Range("C90").Value = "10236900"
$ws.Range("D90").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E90").Value = "1905"
$ws.Range("F90").Value = 43480
$ws.Range("G90").Value = 1087000
$ws.Range("B91").Value = "CC"
$ws.Range("C91").Value = "10236900"
$ws.Range("D91").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E91").Value = "1904"
$ws.Range("F91").Value = 43480
$ws.Range("G91").Value = 1087000
$ws.Range("B92").Value = "CC"
$ws.Range("C92").Value = "10236900"
$ws.Range("D92").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E92").Value = "1903"
$ws.Range("F92").Value = 43480
$ws.Range("G92").Value = 1087000
$ws.Range("B93").Value = "CC"
$ws.Range("C93").Value = "10236900"
$ws.Range("D93").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E93").Value = "1902"
$ws.Range("F93").Value = 43480
$ws.Range("G93").Value = 1087000
$ws.Range("B94").Value = "CC"
$ws.Range("C94").Value = "10236900"
$ws.Range("D94").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E94").Value = "1901"
$ws.Range("F94").Value = 43480
$ws.Range("G94").Value = 1087000
$ws.Range("B95").Value = "CC"
$ws.Range("C95").Value = "10236900"
$ws.Range("D95").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E95").Value = "1812"
$ws.Range("F95").Value = 43480
$ws.Range("G95").Value = 1087000
$ws.Range("B96").Value = "CC"
$ws.Range("C96").Value = "10236900"
$ws.Range("D96").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E96").Value = "1811"
$ws.Range("F96").Value = 43480
$ws.Range("G96").Value = 1087000
$ws.Range("B97").Value = "CC"
$ws.Range("C97").Value = "10236900"
$ws.Range("D97").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E97").Value = "1810"
$ws.Range("F97").Value = 43480
$ws.Range("G97").Value = 1087000
$ws.Range("B98").Value = "CC"
$ws.Range("C98").Value = "10236900"
$ws.Range("D98").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E98").Value = "1809"
$ws.Range("F98").Value = 43480
$ws.Range("G98").Value = 1087000
$ws.Range("B99").Value = "CC"
$ws.Range("C99").Value = "10236900"
$ws.Range("D99").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E99").Value = "1808"
$ws.Range("F99").Value = 43480
$ws.Range("G99").Value = 1087000
$ws.Range("B100").Value = "CC"
$ws.Range("C100").Value = "10236900"
$ws.Range("D100").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E100").Value = "1807"
$ws.Range("F100").Value = 43480
$ws.Range("G100").Value = 1087000
$ws.Range("B101").Value = "CC"
$ws.Range("C101").Value = "10236900"
$ws.Range("D101").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E101").Value = "1806"
$ws.Range("F101").Value = 43480
$ws.Range("G101").Value = 1087000
$ws.Range("B102").Value = "CC"
$ws.Range("C102").Value = "10236900"
$ws.Range("D102").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E102").Value = "1805"
$ws.Range("F102").Value = 43480
$ws.Range("G102").Value = 1087000
$ws.Range("B103").Value = "CC"
$ws.Range("C103").Value = "10236900"
$ws.Range("D103").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E103").Value = "1804"
$ws.Range("F103").Value = 43480
$ws.Range("G103").Value = 1087000
$ws.Range("B104").Value = "CC"
$ws.Range("C104").Value = "10236900"
$ws.Range("D104").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E104").Value = "1803"
$ws.Range("F104").Value = 43480
$ws.Range("G104").Value = 1087000
$ws.Range("B105").Value = "CC"
$ws.Range("C105").Value = "10236900"
$ws.Range("D105").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E105").Value = "1802"
$ws.Range("F105").Value = 43480
$ws.Range("G105").Value = 1087000
$ws.Range("B106").Value = "CC"
$ws.Range("C106").Value = "10236900"
$ws.Range("D106").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E106").Value = "1801"
$ws.Range("F106").Value = 43480
$ws.Range("G106").Value = 1087000
$ws.Range("B107").Value = "CC"
$ws.Range("C107").Value = "10236900"
$ws.Range("D107").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E107").Value = "1712"
$ws.Range("F107").Value = 43480
$ws.Range("G107").Value = 1087000
$ws.Range("B108").Value = "CC"
$ws.Range("C108").Value = "10236900"
$ws.Range("D108").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E108").Value = "1711"
$ws.Range("F108").Value = 43480
$ws.Range("G108").Value = 1087000
$ws.Range("B109").Value = "CC"
$ws.Range("C109").Value = "10236900"
$ws.Range("D109").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E109").Value = "1710"
$ws.Range("F109").Value = 43480
$ws.Range("G109").Value = 1087000
$ws.Range("B110").Value = "CC"
$ws.Range("C110").Value = "10236900"
$ws.Range("D110").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E110").Value = "1709"
$ws.Range("F110").Value = 43480
$ws.Range("G110").Value = 1087000
$ws.Range("B111").Value = "CC"
$ws.Range("C111").Value = "10236900"
$ws.Range("D111").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E111").Value = "1708"
$ws.Range("F111").Value = 43480
$ws.Range("G111").Value = 1087000
$ws.Range("B112").Value = "CC"
$ws.Range("C112").Value = "10236900"
$ws.Range("D112").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E112").Value = "1707"
$ws.Range("F112").Value = 43480
$ws.Range("G112").Value = 1087000
$ws.Range("B113").Value = "CC"
$ws.Range("C113").Value = "10236900"
$ws.Range("D113").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E113").Value = "1706"
$ws.Range("F113").Value = 43480
$ws.Range("G113").Value = 1087000
$ws.Range("B114").Value = "CC"
$ws.Range("C114").Value = "10236900"
$ws.Range("D114").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E114").Value = "1705"
$ws.Range("F114").Value = 43480
$ws.Range("G114").Value = 1087000
$ws.Range("B115").Value = "CC"
$ws.Range("C115").Value = "10236900"
$ws.Range("D115").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E115").Value = "1704"
$ws.Range("F115").Value = 43480
$ws.Range("G115").Value = 1087000
$ws.Range("B116").Value = "CC"
$ws.Range("C116").Value = "10236900"
$ws.Range("D116").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E116").Value = "1703"
$ws.Range("F116").Value = 43480
$ws.Range("G116").Value = 1087000
$ws.Range("B117").Value = "CC"
$ws.Range("C117").Value = "10236900"
$ws.Range("D117").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E117").Value = "1702"
$ws.Range("F117").Value = 43480
$ws.Range("G117").Value = 1087000
$ws.Range("B118").Value = "CC"
$ws.Range("C118").Value = "10236900"
$ws.Range("D118").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E118").Value = "1701"
$ws.Range("F118").Value = 43480
$ws.Range("G118").Value = 1087000
$ws.Range("B119").Value = "CC"
$ws.Range("C119").Value = "10236900"
$ws.Range("D119").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E119").Value = "1612"
$ws.Range("F119").Value = 43480
$ws.Range("G119").Value = 1087000
$ws.Range("B120").Value = "CC"
$ws.Range("C120").Value = "10236900"
$ws.Range("D120").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E120").Value = "1611"
$ws.Range("F120").Value = 43480
$ws.Range("G120").Value = 1087000
$ws.Range("B121").Value = "CC"
$ws.Range("C121").Value = "10236900"
$ws.Range("D121").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E121").Value = "1610"
$ws.Range("F121").Value = 43480
$ws.Range("G121").Value = 1087000
$ws.Range("B122").Value = "CC"
$ws.Range("C122").Value = "10236900"
$ws.Range("D122").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E122").Value = "1609"
$ws.Range("F122").Value = 43480
$ws.Range("G122").Value = 1087000
$ws.Range("B123").Value = "CC"
$ws.Range("C123").Value = "10236900"
$ws.Range("D123").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E123").Value = "1608"
$ws.Range("F123").Value = 43480
$ws.Range("G123").Value = 1087000
$ws.Range("B124").Value = "CC"
$ws.Range("C124").Value = "10236900"
$ws.Range("D124").Value = "OLMEDO DE JESUS CASTRO CAMPEON"
$ws.Range("E124").Value = "1607"
$ws.Range("F124").Value = 43480
$ws.Range("G124").Value = 1087000
$ws.Range("B125").Value = "CC"
$ws.Range("C125").Value = "73561364"
$ws.Range("D125").Value = "RAFAEL ENRIQUE JIMENEZ GONGORA"
$ws.Range("E125").Value = "2507"
$ws.Range("F125").Value = 1898
$ws.Range("G125").Value = 1423500
$ws.Range("B126").Value = "CC"
$ws.Range("C126").Value = "73561364"
$ws.Range("D126").Value = "RAFAEL ENRIQUE JIMENEZ GONGORA"
$ws.Range("E126").Value = "2505"
$ws.Range("F126").Value = 1898
$ws.Range("G126").Value = 1423500
$ws.Range("B127").Value = "CC"
$ws.Range("C127").Value = "73561364"
$ws.Range("D127").Value = "RAFAEL ENRIQUE JIMENEZ GONGORA"
$ws.Range("E127").Value = "2503"
$ws.Range("F127").Value = 1898
$ws.Range("G127").Value = 1423500
$ws.Range("B128").Value = "CC"
$ws.Range("C128").Value = "8802936"
$ws.Range("D128").Value = "CARLOS ALBERTO ARZUZA DIAZ"
$ws.Range("E128").Value = "2003"
$ws.Range("F128").Value = 40892
$ws.Range("G128").Value = 1022300
# Row 128 (last-style row) leaves H/I/J blank, same as the other data rows.
$ws.Range("H128:J128").ClearContents()

Write-Host "edit complete"